$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 81, pushing the old blank separator row (formerly 81)
# and the summary rows (formerly 82-84) down by one.
$ws.Rows.Item(81).Insert() | Out-Null

# Populate the new row 81 with the additional time-tracking entry
# (2014-03-20, 18:15 - 19:00 -> 45 min / 0.75 h).
$ws.Range("A81").Value = 2014
$ws.Range("B81").Value = 3
$ws.Range("C81").Value = 20
$ws.Range("D81").Value = 0.76041666666666663
$ws.Range("E81").Value = 0.79166666666666663
$ws.Range("F81").Formula = "=(E81-D81)*24*60"
$ws.Range("G81").Formula = "=F81/60"

# Reflect the new active cell / selection shown in the GUI.
$ws.Range("F81").Select() | Out-Null
